# Enable data extraction of different studies, using setting study_name in compose
# Update the "Variables" data dictionary sheet:
#  - split icd10_grouped_entities into icd10_grouped + icd10_entity
#  - reorder date_diagnosis to follow the icd10_* columns
#  - add date_diagnosis_year/month/day right after date_diagnosis
#  - tweak several label/description strings
#  - add a trailing gender_mapped row (row 12) after gender (row 11)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Column B = name, C = valueType, D = entityType, L = description
# (A/H/J are identical "df" / 0 / 1 for every data row, before and after)

$rows = @(
    @{ B = "condition_id";          C = "string";  D = "Participant"; L = "Condition ID, unique for each condition" },
    @{ B = "icd10_code";            C = "string";  D = "Participant"; L = "ICD10 GM diagnosis code" },
    @{ B = "icd10_mapped";          C = "decimal"; D = "Participant"; L = "ICD10 GM diagnosis code mapped A = 1, B = 2, C = 3, D = 4,`n        e.g.: A01.9 = 101.9, C50.1 = 350.1 or D41.9 = 441.9" },
    @{ B = "icd10_grouped";         C = "integer"; D = "Participant"; L = "ICD10 GM diagnosis code grouped to parent code, e.g. A01.1`n        and A01.9 both belong to group 101 (remove decimal from icd10_mapped)" },
    @{ B = "icd10_entity";          C = "integer"; D = "Participant"; L = "Entities of resulting ICD10 groups, see utils" },
    @{ B = "date_diagnosis";        C = "string";  D = "Participant"; L = "Date of diagnosis" },
    @{ B = "date_diagnosis_year";   C = "integer"; D = "Participant"; L = "Year of diagnosis" },
    @{ B = "date_diagnosis_month";  C = "integer"; D = "Participant"; L = "Month of Diagnosis" },
    @{ B = "date_diagnosis_day";    C = "integer"; D = "Participant"; L = "Day of Diagnosis" },
    @{ B = "gender";                C = "string";  D = "Participant"; L = "Gender - male, female, other/diverse" },
    @{ B = "gender_mapped";         C = "integer"; D = "Participant"; L = "Gender mapped: 0 = None, 1 = female, 2 = male,`n        3 = other/diverse" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = "df"
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 10).Value = 1
    $ws.Cells.Item($r, 12).Value = $row.L
    $r = $r + 1
}
